$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has the "Description" text sitting in column D
# (under the "Booking.com Price" header) while column E ("Description"
# header) is empty. Add the real Booking.com price into column D and
# push the existing description values over into column E for each
# data row.

$descriptions = @{
    2 = $ws.Range("D2").Value2
    3 = $ws.Range("D3").Value2
    4 = $ws.Range("D4").Value2
    5 = $ws.Range("D5").Value2
    6 = $ws.Range("D6").Value2
}

$prices = @{
    2 = "92.88 USD"
    3 = "391.64 USD"
    4 = "87.82 USD"
    5 = "88.29 USD"
    6 = "109 USD"
}

foreach ($row in 2..6) {
    $ws.Range("E$row").Value = $descriptions[$row]
    $ws.Range("D$row").Value = $prices[$row]
}
